# no-op test
$wb = $excel.ActiveWorkbook
